$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- 1. Refresh the "fetched at" timestamp for every existing data row (2..17) ---
$newDate = "2025-12-04 12:40:02"
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $newDate
}

# --- 2. Push the current row 17 (the mp3-speaker listing) down to a new row 18 ---
$ws.Rows("18:18").Insert()

$ws.Range("A18").Value = $ws.Range("A17").Value2
$ws.Range("B18").Value = $ws.Range("B17").Value2
$ws.Range("C18").Value = $ws.Range("C17").Value2
$ws.Range("D18").Value = $ws.Range("D17").Value2
$ws.Range("E18").Value = $ws.Range("E17").Value2
$ws.Range("F18").Value = $ws.Range("F17").Value2
$ws.Range("G18").Value = $ws.Range("G17").Value2

# --- 3. Overwrite row 17 in place with the newly scraped listing ---
$ws.Range("B17").Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5443568"
$ws.Range("G17").Value = 13

# --- 4. Rebuild the URL hyperlinks (F2:F18) so each points at the right listing ---
$ws.Range("F1").Hyperlinks.Delete()

$urls = @(
  "https://www.lancers.jp/work/detail/5423720",
  "https://www.lancers.jp/work/detail/5446833",
  "https://www.lancers.jp/work/detail/5419380",
  "https://www.lancers.jp/work/detail/5447137",
  "https://www.lancers.jp/work/detail/5446990",
  "https://www.lancers.jp/work/detail/5446867",
  "https://www.lancers.jp/work/detail/5441557",
  "https://www.lancers.jp/work/detail/5447021",
  "https://www.lancers.jp/work/detail/5441568",
  "https://www.lancers.jp/work/detail/5446668",
  "https://www.lancers.jp/work/detail/5431107",
  "https://www.lancers.jp/work/detail/5446849",
  "https://www.lancers.jp/work/detail/5446997",
  "https://www.lancers.jp/work/detail/5447102",
  "https://www.lancers.jp/work/detail/5437544",
  "https://www.lancers.jp/work/detail/5443568",
  "https://www.lancers.jp/work/detail/5446806"
)

for ($i = 0; $i -lt $urls.Count; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $urls[$i]) | Out-Null
}

# --- 5. Sheet dimension now spans one more row ---
$ws.Range("A1:H18").Select()
